# Append: 2025-11-19 12:49 JST
# Update the "取得日時" (acquired timestamp) column (A) for rows 2-11 on the
# "ランサーズ" sheet from the previous run's timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-19 12:49:55"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
